# BOM finished, Sourcing merged
# Add the new "PIFACE CONTROL & DISPLAY 2" line item to the Bill of Materials,
# including its hyperlink, and leave a formatted (but empty) cost cell on the
# following row ready for the next entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 8: PiFace Control & Display 2 -------------------------------
$ws.Range("A8").Value = "PIFACE CONTROL & DISPLAY 2"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 26.49
$ws.Range("D8").Value = 2434231

$link = "https://fi.farnell.com/piface/piface-control-display-2/i-o-board-w-lcd-for-raspberry/dp/2434231?st=raspberry%20display"
$ws.Range("E8").Value = $link

# Currency format (2 decimals) for the new cost cell
$ws.Range("C8").NumberFormat = "#,##0.00\ [$€-1];[Red]\-#,##0.00\ [$€-1]"

# Prime the next row's cost cell with the normal currency format (no decimals),
# same as the rest of the table, even though it is otherwise still empty.
$ws.Range("C9").NumberFormat = "#,##0\ [$€-1];[Red]\-#,##0\ [$€-1]"

# Turn the new link cell into a real hyperlink
$ws.Hyperlinks.Add($ws.Range("E8"), $link)

# Update dimension / selection so the workbook re-opens focused on the new row
$null = $ws.Activate()
$null = $ws.Range("A9:E9").Select()
